$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell (well outside the used range) used to stage values that Excel
# would otherwise auto-convert from text to a number (e.g. "0.17", "-0.01").
# We format it as Text, assign the value, copy it, and paste-special just the
# value into the destination cell. Because PasteSpecial(xlPasteValues) only
# transfers the value (and its string/number type), the destination cell
# keeps its original (default) style while still being stored as a shared
# string of the exact text. The helper cell is fully cleared (Clear, not
# just ClearContents) afterwards so it leaves no trace in the sheet data or
# used range.
$xlPasteValues = -4163
$helper = $ws.Range("Z100")

function Set-TextValue($cellAddress, $text) {
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($cellAddress).PasteSpecial($xlPasteValues)
}

# Row 2 (A Lag)
Set-TextValue "B2" "0.17"
$ws.Range("C2").Value = "44.29***"
Set-TextValue "D2" "-0.89"

# Row 3 (FFR Lag)
Set-TextValue "B3" "-0.01"
$ws.Range("C3").Value = "2.21***"
$ws.Range("D3").Value = "0.46***"

# Row 4 (LF Lag)
Set-TextValue "B4" "-0.09"
Set-TextValue "C4" "0.98"
$ws.Range("D4").Value = "0.82*"

# Clean up the helper cell completely (contents + formatting) so it does not
# appear in the saved worksheet.
$helper.Clear()
